$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new blank column before column N ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$null = $wsRepay.Columns("N:N").Insert()
# match the width of the neighbouring "In Advance" column (M, width 11)
$wsRepay.Columns("N:N").ColumnWidth = 10.166666666666666

# --- Transactions sheet keeps its own saved selection, but is no longer active ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$null = $wsTrans.Range("B9").Select()

# Make "Repayment schedule" the active sheet and update its selection (done last
# so it ends up being the active/selected tab when the workbook is saved)
$null = $wsRepay.Activate()
$null = $wsRepay.Range("J18").Select()
